$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("C4").Value = 4367
$ws.Range("D4").Value = 82.8

# Row 5
$ws.Range("C5").Value = 295
$ws.Range("D5").Value = 96.5

# Row 6
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 848
$ws.Range("D6").Value = 96.59999999999999

# Row 7 - driver text swaps with row 8
$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.60.1.2"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 464
$ws.Range("D7").Value = 98.09999999999999

# Row 8
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.2.1"
$ws.Range("B8").Value = 44
$ws.Range("C8").Value = 1791
$ws.Range("D8").Value = 98.2

# Row 9
$ws.Range("B9").Value = 250
$ws.Range("C9").Value = 7654
$ws.Range("D9").Value = 98.3

# Row 10
$ws.Range("B10").Value = 511
$ws.Range("C10").Value = 18859

# Row 11
$ws.Range("B11").Value = 293
$ws.Range("C11").Value = 5870

# Row 12 - Totals
$ws.Range("B12").Value = 1127
$ws.Range("C12").Value = 40175

# Row 20
$ws.Range("B20").Value = 449371

# Row 24
$ws.Range("B24").Value = 77999
